# Add data organization files for ESD
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("samples_retained")

# Row 16: rename the old combined "esd" row to "esd [en]" and fill in its
# sample counts / language, keeping the note column text the same.
$ws.Range("A16").Value = "esd [en]"
$ws.Range("B16").Value = "acted"
$ws.Range("C16").Value = 3500
$ws.Range("D16").Value = 10500
$ws.Range("E16").Value = 3500
$ws.Range("F16").Value = "English"
$ws.Range("G16").Formula = '=IF(OR(ISBLANK(C16), ISBLANK(D16),ISBLANK(E16)), "", SUM(C16:E16))'
$ws.Range("H16").Value = "English and Mandarin"

# Row 17: new row for the Mandarin portion of ESD.
$ws.Range("A17").Value = "esd [zh]"
$ws.Range("B17").Value = "acted"
$ws.Range("C17").Value = 3500
$ws.Range("D17").Value = 10500
$ws.Range("E17").Value = 3500
$ws.Range("F17").Value = "Mandarin Chinese"
$ws.Range("G17").Formula = '=IF(OR(ISBLANK(C17), ISBLANK(D17),ISBLANK(E17)), "", SUM(C17:E17))'
$ws.Range("H17").Value = "English and Mandarin"

# Widen column G slightly to fit the new values (Excel stores column widths
# with a small fixed offset from the COM ColumnWidth value, so back that
# offset out to land exactly on a displayed width of 6) and move the
# selection like the author's saved view.
$ws.Columns.Item(7).ColumnWidth = 5.166666666666667
$ws.Range("J17").Select()
